# Auto-generated financial data update for SHI worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SHI")

$ws.Range("D8").Value = 11756800
$ws.Range("E8").Value = 9785600
$ws.Range("F8").Value = 9949000
$ws.Range("G8").Value = 13761300
$ws.Range("H8").Value = 15657700
$ws.Range("I8").Value = 12943900
$ws.Range("J8").Value = 13284100
$ws.Range("D9").Value = 10744600
$ws.Range("E9").Value = 8716400
$ws.Range("F9").Value = 9313800
$ws.Range("G9").Value = 13788800
$ws.Range("H9").Value = 15319800
$ws.Range("I9").Value = 13151800
$ws.Range("J9").Value = 13042400
$ws.Range("D10").Value = 1012200
$ws.Range("E10").Value = 1069300
$ws.Range("F10").Value = 635200
$ws.Range("G10").Value = -27500
$ws.Range("H10").Value = 338000
$ws.Range("I10").Value = -207900
$ws.Range("J10").Value = 241700
$ws.Range("G14").Value = -1000
$ws.Range("E15").Value = 1300
$ws.Range("D17").Value = 10806700
$ws.Range("E17").Value = 8779700
$ws.Range("F17").Value = 9368900
$ws.Range("G17").Value = 13848600
$ws.Range("H17").Value = 15332400
$ws.Range("I17").Value = 13207000
$ws.Range("J17").Value = 13126800
$ws.Range("D18").Value = 950100
$ws.Range("E18").Value = 1005900
$ws.Range("F18").Value = 580100
$ws.Range("G18").Value = -87300
$ws.Range("H18").Value = 325400
$ws.Range("I18").Value = -263000
$ws.Range("J18").Value = 157300
$ws.Range("D20").Value = 223400
$ws.Range("E20").Value = 156400
$ws.Range("F20").Value = 80200
$ws.Range("G20").Value = 10800
$ws.Range("H20").Value = 37500
$ws.Range("I20").Value = 16600
$ws.Range("J20").Value = 67100
$ws.Range("D21").Value = 1410100
$ws.Range("E21").Value = 1406300
$ws.Range("F21").Value = 916400
$ws.Range("H21").Value = 742200
$ws.Range("D22").Value = 8100
$ws.Range("E22").Value = 8000
$ws.Range("F22").Value = 31500
$ws.Range("G22").Value = 55600
$ws.Range("I22").Value = 52800
$ws.Range("J22").Value = 32000
$ws.Range("D23").Value = 1165400
$ws.Range("E23").Value = 1154400
$ws.Range("F23").Value = 628800
$ws.Range("G23").Value = -132100
$ws.Range("H23").Value = 362800
$ws.Range("I23").Value = -299300
$ws.Range("J23").Value = 192400
$ws.Range("D24").Value = 252100
$ws.Range("E24").Value = 266700
$ws.Range("F24").Value = 137500
$ws.Range("G24").Value = -31800
$ws.Range("H24").Value = 56300
$ws.Range("I24").Value = -75900
$ws.Range("J24").Value = 46000
$ws.Range("D26").Value = 913300
$ws.Range("E26").Value = 887700
$ws.Range("F26").Value = 491300
$ws.Range("G26").Value = -100300
$ws.Range("H26").Value = 306500
$ws.Range("I26").Value = -223400
$ws.Range("J26").Value = 146400
$ws.Range("D27").Value = 911700
$ws.Range("E27").Value = 885800
$ws.Range("F27").Value = 485900
$ws.Range("G27").Value = -102700
$ws.Range("H27").Value = 305000
$ws.Range("I27").Value = -226800
$ws.Range("J27").Value = 141900
$ws.Range("D32").Value = -223400
$ws.Range("E32").Value = -156400
$ws.Range("F32").Value = -80200
$ws.Range("G32").Value = -10800
$ws.Range("H32").Value = -37500
$ws.Range("I32").Value = -16600
$ws.Range("J32").Value = -67100
$ws.Range("D33").Value = 911700
$ws.Range("E33").Value = 885800
$ws.Range("F33").Value = 485900
$ws.Range("G33").Value = -102700
$ws.Range("H33").Value = 305000
$ws.Range("I33").Value = -226800
$ws.Range("J33").Value = 141900
$ws.Range("D35").Value = 911700
$ws.Range("E35").Value = 885800
$ws.Range("F35").Value = 485900
$ws.Range("G35").Value = -102700
$ws.Range("H35").Value = 305000
$ws.Range("I35").Value = -226800
$ws.Range("J35").Value = 141900
$ws.Range("D41").Value = 1109400
$ws.Range("E41").Value = 782300
$ws.Range("F41").Value = 159100
$ws.Range("G41").Value = 40700
$ws.Range("H41").Value = 18700
$ws.Range("I41").Value = 23900
$ws.Range("J41").Value = 13600
$ws.Range("D42").Value = 301100
$ws.Range("E42").Value = 25100
$ws.Range("H42").Value = 1100
$ws.Range("D43").Value = 524800
$ws.Range("E43").Value = 461600
$ws.Range("F43").Value = 428700
$ws.Range("G43").Value = 489800
$ws.Range("H43").Value = 788600
$ws.Range("I43").Value = 866600
$ws.Range("J43").Value = 592500
$ws.Range("D44").Value = 979100
$ws.Range("E44").Value = 914100
$ws.Range("F44").Value = 620100
$ws.Range("G44").Value = 880200
$ws.Range("H44").Value = 1341500
$ws.Range("I44").Value = 1385100
$ws.Range("J44").Value = 828500
$ws.Range("D45").Value = 33900
$ws.Range("E45").Value = 24600
$ws.Range("D46").Value = 2948300
$ws.Range("E46").Value = 2207700
$ws.Range("F46").Value = 1208600
$ws.Range("G46").Value = 1411400
$ws.Range("H46").Value = 2149900
$ws.Range("I46").Value = 1913200
$ws.Range("J46").Value = 1434500
$ws.Range("D47").Value = 660700
$ws.Range("E47").Value = 547500
$ws.Range("F47").Value = 491400
$ws.Range("G47").Value = 435800
$ws.Range("H47").Value = 444300
$ws.Range("I47").Value = 425500
$ws.Range("J47").Value = 430600
$ws.Range("D48").Value = 2116200
$ws.Range("E48").Value = 2162700
$ws.Range("F48").Value = 2302000
$ws.Range("G48").Value = 2448800
$ws.Range("H48").Value = 2605400
$ws.Range("I48").Value = 5406300
$ws.Range("J48").Value = 2494400
$ws.Range("D49").Value = 55400
$ws.Range("E49").Value = 48300
$ws.Range("F49").Value = 57700
$ws.Range("G49").Value = 94200
$ws.Range("H49").Value = 73300
$ws.Range("I49").Value = 267500
$ws.Range("J49").Value = 122500
$ws.Range("D52").Value = 73200
$ws.Range("E52").Value = 71700
$ws.Range("F52").Value = 69000
$ws.Range("G52").Value = 196500
$ws.Range("H52").Value = 164400
$ws.Range("I52").Value = 380600
$ws.Range("J52").Value = 77100
$ws.Range("D54").Value = 5853800
$ws.Range("E54").Value = 5037900
$ws.Range("F54").Value = 4128900
$ws.Range("G54").Value = 4586700
$ws.Range("H54").Value = 5437300
$ws.Range("I54").Value = 5411400
$ws.Range("J54").Value = 4559000
$ws.Range("D57").Value = 283200
$ws.Range("E57").Value = 315200
$ws.Range("F57").Value = 231900
$ws.Range("G57").Value = 521100
$ws.Range("H57").Value = 406600
$ws.Range("I57").Value = 856800
$ws.Range("J57").Value = 464000
$ws.Range("D58").Value = 90000
$ws.Range("E58").Value = 81100
$ws.Range("F58").Value = 307200
$ws.Range("G58").Value = 605200
$ws.Range("H58").Value = 1052800
$ws.Range("I58").Value = 3272100
$ws.Range("J58").Value = 818000
$ws.Range("D59").Value = 1247800
$ws.Range("E59").Value = 930800
$ws.Range("F59").Value = 607600
$ws.Range("G59").Value = 726500
$ws.Range("H59").Value = 1214500
$ws.Range("I59").Value = 982400
$ws.Range("J59").Value = 539200
$ws.Range("D60").Value = 1621000
$ws.Range("E60").Value = 1327100
$ws.Range("F60").Value = 1146700
$ws.Range("G60").Value = 1852900
$ws.Range("H60").Value = 2674000
$ws.Range("I60").Value = 2809000
$ws.Range("J60").Value = 1821300
$ws.Range("G61").Value = 242300
$ws.Range("H61").Value = 93200
$ws.Range("I61").Value = 182700
$ws.Range("J61").Value = 23800
$ws.Range("J62").Value = 13600
$ws.Range("D66").Value = 1664200
$ws.Range("E66").Value = 1368900
$ws.Range("F66").Value = 1190700
$ws.Range("G66").Value = 2137900
$ws.Range("H66").Value = 2805600
$ws.Range("I66").Value = 3031300
$ws.Range("J66").Value = 1898700
$ws.Range("D72").Value = 2575500
$ws.Range("E72").Value = 2066200
$ws.Range("F72").Value = 1335300
$ws.Range("G72").Value = 846000
$ws.Range("H72").Value = 1028900
$ws.Range("I72").Value = 1903300
$ws.Range("J72").Value = 1232500
$ws.Range("D76").Value = 4189600
$ws.Range("E76").Value = 3669000
$ws.Range("F76").Value = 2938100
$ws.Range("G76").Value = 2448800
$ws.Range("H76").Value = 2631700
$ws.Range("I76").Value = 2380100
$ws.Range("J76").Value = 2660300
$ws.Range("D81").Value = 911700
$ws.Range("E81").Value = 885800
$ws.Range("F81").Value = 485900
$ws.Range("G81").Value = -102700
$ws.Range("H81").Value = 305000
$ws.Range("I81").Value = -226800
$ws.Range("J81").Value = 141900
$ws.Range("D83").Value = 236200
$ws.Range("E83").Value = 243500
$ws.Range("F83").Value = 255600
$ws.Range("H83").Value = 378700
$ws.Range("D89").Value = 1047900
$ws.Range("E89").Value = 1065800
$ws.Range("F89").Value = 732100
$ws.Range("G89").Value = 543500
$ws.Range("H89").Value = 756700
$ws.Range("I89").Value = -306700
$ws.Range("J89").Value = 329500
$ws.Range("D91").Value = -177700
$ws.Range("E91").Value = -133800
$ws.Range("F91").Value = -103200
$ws.Range("G91").Value = -161700
$ws.Range("H91").Value = -196400
$ws.Range("I91").Value = -632200
$ws.Range("J91").Value = -516700
$ws.Range("D94").Value = -356300
$ws.Range("E94").Value = -28200
$ws.Range("F94").Value = -65100
$ws.Range("G94").Value = -135100
$ws.Range("H94").Value = -93400
$ws.Range("I94").Value = -602900
$ws.Range("J94").Value = "NA"
$ws.Range("D96").Value = -400300
$ws.Range("E96").Value = -161000
$ws.Range("G96").Value = -80600
$ws.Range("H96").Value = -53500
$ws.Range("I96").Value = -53600
$ws.Range("J96").Value = -105800
$ws.Range("D100").Value = -384300
$ws.Range("E100").Value = -391400
$ws.Range("F100").Value = -548500
$ws.Range("G100").Value = -386800
$ws.Range("H100").Value = -667400
$ws.Range("I100").Value = 919900
$ws.Range("J100").Value = "NA"
$ws.Range("E101").Value = 1300
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = 306300
$ws.Range("E102").Value = 647500
$ws.Range("F102").Value = 118500
$ws.Range("G102").Value = 21700
$ws.Range("H102").Value = -4100
$ws.Range("I102").Value = 10300
